$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.05
$ws.Range("F2").Value = 0
$ws.Range("H2").Value = 0.08
$ws.Range("I2").Value = -0.08
$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 0.13
$ws.Range("F3").Value = 0.01
$ws.Range("H3").Value = 0.14
$ws.Range("I3").Value = -0.14
$ws.Range("B4").Value = 0
$ws.Range("D4").Value = 0.06
$ws.Range("F4").Value = 0
$ws.Range("I4").Value = -0.09
$ws.Range("B5").Value = 0
$ws.Range("D5").Value = 0.04
$ws.Range("E5").Value = -0.03
$ws.Range("F5").Value = 0
$ws.Range("H5").Value = 0.06
$ws.Range("I5").Value = -0.06
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0.05
$ws.Range("H6").Value = 0.12
$ws.Range("I6").Value = -0.08
$ws.Range("D7").Value = 0.06
$ws.Range("E7").Value = -0.06
$ws.Range("F7").Value = 0
$ws.Range("H7").Value = 0.09
$ws.Range("I7").Value = -0.06
$ws.Range("B8").Value = 0
$ws.Range("E8").Value = -0.05
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0.04
$ws.Range("B9").Value = 0.01
$ws.Range("D9").Value = 0.08
$ws.Range("E9").Value = -0.05
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0.05
$ws.Range("H9").Value = 0.14
$ws.Range("I9").Value = -0.11
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0.02
$ws.Range("D10").Value = 0.04
$ws.Range("G10").Value = 0.03
$ws.Range("H10").Value = 0.06
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0.03
$ws.Range("E11").Value = -0.06
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0.04
$ws.Range("H11").Value = 0.1
$ws.Range("I11").Value = -0.08
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0.03
$ws.Range("D12").Value = 0.06
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0.04
$ws.Range("H12").Value = 0.1
$ws.Range("I12").Value = -0.08
$ws.Range("B13").Value = 0
$ws.Range("D13").Value = 0.05
$ws.Range("E13").Value = -0.06
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0.04
$ws.Range("H13").Value = 0.08
$ws.Range("I13").Value = -0.09
$ws.Range("E14").Value = -0.03
$ws.Range("F14").Value = 0
$ws.Range("H14").Value = 0.06
$ws.Range("I14").Value = -0.04
$ws.Range("F15").Value = 0
$ws.Range("H15").Value = 0.07000000000000001
$ws.Range("I15").Value = -0.07000000000000001
$ws.Range("D16").Value = 0.04
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0.03
$ws.Range("I16").Value = -0.07000000000000001
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0.02
$ws.Range("D17").Value = 0.04
$ws.Range("E17").Value = -0.04
$ws.Range("F17").Value = 0
$ws.Range("H17").Value = 0.08
$ws.Range("I17").Value = -0.06
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0.04
$ws.Range("D18").Value = 0.11
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0.05
$ws.Range("H18").Value = 0.14
$ws.Range("I18").Value = -0.12
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0.03
$ws.Range("D19").Value = 0.08
$ws.Range("F19").Value = -0
$ws.Range("H19").Value = 0.13
$ws.Range("E20").Value = -0.1
$ws.Range("F20").Value = -0
$ws.Range("G20").Value = 0.06
$ws.Range("H20").Value = 0.14
$ws.Range("I20").Value = -0.14
$ws.Range("D21").Value = 0.09
$ws.Range("E21").Value = -0.09
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0.04
$ws.Range("H21").Value = 0.14
$ws.Range("I21").Value = -0.1
